# Natmi following Dr Hou advice
# Update recomputed NATMI ligand-receptor statistics for Fbn1-Itgav pairs.
# Ligand-expressing cells (E) and Receptor-expressing cells (K) changed
# from 1 to 3, and all the downstream derived columns (G,H,I,J,M,N,O,P,Q,R,S,T)
# were recalculated accordingly for data rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; G=7.871715666666666;  H=23.615147;          I=0.02771913691218268; J=0.02771913691218268; K=3; M=21.09934133333334; N=63.29802400000001;  O=0.2917236204149438; P=0.2917236204149438; Q=166.0880157299476;  R=1494.792141569528;  S=0.008086326974799437; T=0.008086326974799437 }
    3  = @{ E=3; G=7.871715666666666;  H=23.615147;          I=0.02771913691218268; J=0.02771913691218268; K=3; M=35.81943766666667; N=107.458313;         O=0.4952465516465762; P=0.4952465516465762; Q=281.9604286518901;  R=2537.643857867011;  S=0.0137278069703778;   T=0.0137278069703778 }
    4  = @{ E=3; G=7.871715666666666;  H=23.615147;          I=0.02771913691218268; J=0.02771913691218268; K=3; M=15.40769666666667; N=46.22309;           O=0.2130298279384801; P=0.2130298279384801; Q=121.2850072382478;  R=1091.56506514423;   S=0.005905002967005448; T=0.005905002967005448 }
    5  = @{ E=3; G=266.1315866666666;  H=798.3947599999999; I=0.9371448614065047; J=0.9371448614065045; K=3; M=21.09934133333334; N=63.29802400000001;  O=0.2917236204149438; P=0.2917236204149438; Q=5615.201186661582;  R=50536.81067995424;  S=0.2733872918227662;   T=0.2733872918227663 }
    6  = @{ E=3; G=266.1315866666666;  H=798.3947599999999; I=0.9371448614065047; J=0.9371448614065045; K=3; M=35.81943766666667; N=107.458313;         O=0.4952465516465762; P=0.4952465516465762; Q=9532.683779737763;  R=85794.15401763987;  S=0.4641177610048799;   T=0.4641177610048799 }
    7  = @{ E=3; G=266.1315866666666;  H=798.3947599999999; I=0.9371448614065047; J=0.9371448614065045; K=3; M=15.40769666666667; N=46.22309;           O=0.2130298279384801; P=0.2130298279384801; Q=4100.47476077871;   R=36904.2728470084;   S=0.1996398085788584;   T=0.1996398085788584 }
    8  = @{ E=3; G=9.977966333333333;  H=29.933899;          I=0.03513600168131278; J=0.03513600168131277; K=3; M=21.09934133333334; N=63.29802400000001;  O=0.2917236204149438; P=0.2917236204149438; Q=210.5285174795085;  R=1894.756657315576;  S=0.01025000161737811;  T=0.01025000161737811 }
    9  = @{ E=3; G=9.977966333333333;  H=29.933899;          I=0.03513600168131278; J=0.03513600168131277; K=3; M=35.81943766666667; N=107.458313;         O=0.4952465516465762; P=0.4952465516465762; Q=357.4051431169319;  R=3216.646288052387;  S=0.01740098367131845;  T=0.01740098367131845 }
    10 = @{ E=3; G=9.977966333333333;  H=29.933899;          I=0.03513600168131278; J=0.03513600168131277; K=3; M=15.40769666666667; N=46.22309;           O=0.2130298279384801; P=0.2130298279384801; Q=153.7374786142122;  R=1383.63730752791;   S=0.007485016392616207; T=0.007485016392616207 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}
